$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 42; $row++) {
    $cell = $ws.Cells.Item($row, 6)  # Column F = 6
    $cell.Value = $cell.Value2 - 120
}
